{"js": "// Highlight (yellow) the \"Le compte administrateur peut :\" paragraph and the\n// three bullet items that follow it (\"Ajouter des utilisateurs;\",\n// \"Lister les utilisateurs;\", \"Effacer des utilisateurs;\").\n// Matches commit \"ajout autorisation cr\u00e9ation de code\": the admin-account\n// permissions block gets a yellow highlight to call it out.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Exact (trimmed) paragraph texts targeted by the edit, in document order.\nconst targets = [\n  \"Le compte administrateur peut :\",\n  \"Ajouter des utilisateurs;\",\n  \"Lister les utilisateurs;\",\n  \"Effacer des utilisateurs;\"\n];\n\nfor (const paragraph of paragraphs.items) {\n  const text = (paragraph.text || \"\").trim();\n  if (targets.includes(text)) {\n    // Setting highlightColor on the paragraph's font colors every run in the\n    // paragraph AND the paragraph mark itself (pPr/rPr), matching the diff.\n    paragraph.font.highlightColor = \"Yellow\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Highlight (yellow) the \"Le compte administrateur peut :\" paragraph and the\n# three bullet items that follow it (\"Ajouter des utilisateurs;\",\n# \"Lister les utilisateurs;\", \"Effacer des utilisateurs;\").\n# Matches commit \"ajout autorisation cr\u00e9ation de code\": the admin-account\n# permissions block gets a yellow highlight to call it out.\n\n$d = $word.ActiveDocument\n\n# Exact (trimmed) paragraph texts targeted by the edit, in document order.\n$targets = @(\n  \"Le compte administrateur peut :\",\n  \"Ajouter des utilisateurs;\",\n  \"Lister les utilisateurs;\",\n  \"Effacer des utilisateurs;\"\n)\n\nforeach ($p in $d.Paragraphs) {\n  $text = $p.Range.Text.Trim()\n  if ($targets -contains $text) {\n    # Paragraph-level Font highlights every run in the paragraph AND the\n    # paragraph mark itself (pPr/rPr), matching the diff.\n    $p.Range.Font.HighlightColorIndex = 7\n  }\n}\n"}
